{"js": "// Replace the 25 division-problem text values in the practice-sheet table.\n// Old values are each unique within the document, so a targeted search +\n// insertText(\"Replace\") for each pair reproduces the diff exactly while\n// leaving run formatting (rFonts, sz, etc.) untouched.\nconst replacements = [\n  [\"781\u00f73=\", \"721\u00f79=\"],\n  [\"737\u00f77=\", \"258\u00f74=\"],\n  [\"663\u00f79=\", \"996\u00f72=\"],\n  [\"441\u00f79=\", \"110\u00f73=\"],\n  [\"469\u00f77=\", \"149\u00f73=\"],\n  [\"708\u00f77=\", \"974\u00f78=\"],\n  [\"862\u00f77=\", \"387\u00f77=\"],\n  [\"173\u00f76=\", \"178\u00f75=\"],\n  [\"271\u00f78=\", \"268\u00f76=\"],\n  [\"567\u00f77=\", \"461\u00f72=\"],\n  [\"360\u00f73=\", \"692\u00f75=\"],\n  [\"835\u00f79=\", \"377\u00f72=\"],\n  [\"833\u00f79=\", \"101\u00f74=\"],\n  [\"749\u00f78=\", \"385\u00f78=\"],\n  [\"105\u00f79=\", \"759\u00f77=\"],\n  [\"507\u00f76=\", \"773\u00f77=\"],\n  [\"384\u00f72=\", \"101\u00f72=\"],\n  [\"554\u00f74=\", \"583\u00f79=\"],\n  [\"860\u00f74=\", \"558\u00f73=\"],\n  [\"484\u00f74=\", \"979\u00f74=\"],\n  [\"794\u00f76=\", \"467\u00f76=\"],\n  [\"716\u00f76=\", \"926\u00f75=\"],\n  [\"712\u00f73=\", \"554\u00f78=\"],\n  [\"769\u00f73=\", \"968\u00f74=\"],\n  [\"179\u00f78=\", \"405\u00f79=\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 division-problem text values in the practice-sheet table.\n# Each old value is unique within the document, so a targeted Find/Replace\n# for each pair reproduces the diff exactly while leaving run formatting\n# (rFonts, sz, etc.) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"781\u00f73=\"; New = \"721\u00f79=\" },\n    @{ Old = \"737\u00f77=\"; New = \"258\u00f74=\" },\n    @{ Old = \"663\u00f79=\"; New = \"996\u00f72=\" },\n    @{ Old = \"441\u00f79=\"; New = \"110\u00f73=\" },\n    @{ Old = \"469\u00f77=\"; New = \"149\u00f73=\" },\n    @{ Old = \"708\u00f77=\"; New = \"974\u00f78=\" },\n    @{ Old = \"862\u00f77=\"; New = \"387\u00f77=\" },\n    @{ Old = \"173\u00f76=\"; New = \"178\u00f75=\" },\n    @{ Old = \"271\u00f78=\"; New = \"268\u00f76=\" },\n    @{ Old = \"567\u00f77=\"; New = \"461\u00f72=\" },\n    @{ Old = \"360\u00f73=\"; New = \"692\u00f75=\" },\n    @{ Old = \"835\u00f79=\"; New = \"377\u00f72=\" },\n    @{ Old = \"833\u00f79=\"; New = \"101\u00f74=\" },\n    @{ Old = \"749\u00f78=\"; New = \"385\u00f78=\" },\n    @{ Old = \"105\u00f79=\"; New = \"759\u00f77=\" },\n    @{ Old = \"507\u00f76=\"; New = \"773\u00f77=\" },\n    @{ Old = \"384\u00f72=\"; New = \"101\u00f72=\" },\n    @{ Old = \"554\u00f74=\"; New = \"583\u00f79=\" },\n    @{ Old = \"860\u00f74=\"; New = \"558\u00f73=\" },\n    @{ Old = \"484\u00f74=\"; New = \"979\u00f74=\" },\n    @{ Old = \"794\u00f76=\"; New = \"467\u00f76=\" },\n    @{ Old = \"716\u00f76=\"; New = \"926\u00f75=\" },\n    @{ Old = \"712\u00f73=\"; New = \"554\u00f78=\" },\n    @{ Old = \"769\u00f73=\"; New = \"968\u00f74=\" },\n    @{ Old = \"179\u00f78=\"; New = \"405\u00f79=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
